# Scheduled runner update: refresh Universalis market-price snapshots
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / profit
# columns) across the ALC, ARM, BSM, CRP, CUL, GSM and LTW crafting
# sheets in the Maduin_Profits workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value2 = 131.95
$ws.Range("I33").Value2 = 156.45454
$ws.Range("J33").Value2 = 102
$ws.Range("K33").Value2 = 156.45454
$ws.Range("L33").Value2 = 102
$ws.Range("M33").Value2 = 72.54545999999999
$ws.Range("N33").Value2 = -560

$ws.Range("H53").Value2 = 510.33334
$ws.Range("I53").Value2 = 299.66666
$ws.Range("J53").Value2 = 791.2222
$ws.Range("K53").Value2 = 299.66666
$ws.Range("L53").Value2 = 791.2222
$ws.Range("M53").Value2 = 337.33334
$ws.Range("N53").Value2 = -2065.2222

$ws.Range("H86").Value2 = 3903.0833
$ws.Range("I86").Value2 = 3804.111
$ws.Range("K86").Value2 = 3804.111
$ws.Range("M86").Value2 = -2681.111

$ws.Range("H89").Value2 = 3903.0833
$ws.Range("I89").Value2 = 3804.111
$ws.Range("K89").Value2 = 19020.555
$ws.Range("M89").Value2 = -13404.555

$ws.Range("H92").Value2 = 451.69232
$ws.Range("I92").Value2 = 543.3
$ws.Range("K92").Value2 = 543.3
$ws.Range("M92").Value2 = 704.7

$ws.Range("H106").Value2 = 3679.8
$ws.Range("I106").Value2 = 3679.8
$ws.Range("K106").Value2 = 3679.8
$ws.Range("M106").Value2 = -3048.8

$ws.Range("H123").Value2 = 73999
$ws.Range("I123").Value2 = 73999
$ws.Range("K123").Value2 = 73999
$ws.Range("M123").Value2 = -69099

$ws.Range("H132").Value2 = 5499.75
$ws.Range("J132").Value2 = 5000
$ws.Range("L132").Value2 = 15000
$ws.Range("N132").Value2 = -20060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 3206
$ws.Range("I32").Value2 = 3206
$ws.Range("K32").Value2 = 3206
$ws.Range("M32").Value2 = -2919

$ws.Range("H45").Value2 = 5172.2
$ws.Range("I45").Value2 = 1847
$ws.Range("K45").Value2 = 1847
$ws.Range("M45").Value2 = -1470

$ws.Range("H61").Value2 = 3806.7273
$ws.Range("I61").Value2 = 3426.4119
$ws.Range("J61").Value2 = 5099.8
$ws.Range("K61").Value2 = 3426.4119
$ws.Range("L61").Value2 = 5099.8
$ws.Range("M61").Value2 = -3214.4119
$ws.Range("N61").Value2 = -5523.8

$ws.Range("H97").Value2 = 75.916664
$ws.Range("I97").Value2 = 95.5
$ws.Range("K97").Value2 = 95.5
$ws.Range("M97").Value2 = 400.5

$ws.Range("H102").Value2 = 2149.9167
$ws.Range("I102").Value2 = 685.7143
$ws.Range("K102").Value2 = 685.7143
$ws.Range("M102").Value2 = 936.2857

$ws.Range("H122").Value2 = 1832.625
$ws.Range("I122").Value2 = 2315.6667
$ws.Range("K122").Value2 = 6947.000100000001
$ws.Range("M122").Value2 = -4497.000100000001

$ws.Range("H136").Value2 = 3806.7273
$ws.Range("I136").Value2 = 3426.4119
$ws.Range("J136").Value2 = 5099.8
$ws.Range("K136").Value2 = 10279.2357
$ws.Range("L136").Value2 = 15299.4
$ws.Range("M136").Value2 = -7729.235700000001
$ws.Range("N136").Value2 = -20399.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 1753.8823
$ws.Range("I20").Value2 = 1687.7
$ws.Range("J20").Value2 = 1848.4286
$ws.Range("K20").Value2 = 1687.7
$ws.Range("L20").Value2 = 1848.4286
$ws.Range("M20").Value2 = -1440.7
$ws.Range("N20").Value2 = -2342.4286

$ws.Range("H94").Value2 = 2504.7273
$ws.Range("I94").Value2 = 425.33334
$ws.Range("K94").Value2 = 425.33334
$ws.Range("M94").Value2 = 25.66665999999998

$ws.Range("H99").Value2 = 1738.9286
$ws.Range("I99").Value2 = 1558.6364
$ws.Range("K99").Value2 = 1558.6364
$ws.Range("M99").Value2 = -60.63640000000009

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 892.375
$ws.Range("I16").Value2 = 892.375
$ws.Range("J16").Value2 = 0
$ws.Range("K16").Value2 = 892.375
$ws.Range("L16").Value2 = 0
$ws.Range("M16").Value2 = -605.375
$ws.Range("N16").ClearContents()

$ws.Range("H39").Value2 = 0
$ws.Range("I39").Value2 = 0
$ws.Range("K39").Value2 = 0
$ws.Range("M39").ClearContents()

$ws.Range("H43").Value2 = 21885.666
$ws.Range("J43").Value2 = 21885.666
$ws.Range("L43").Value2 = 21885.666
$ws.Range("N43").Value2 = -22253.666

$ws.Range("H49").Value2 = 0
$ws.Range("I49").Value2 = 0
$ws.Range("K49").Value2 = 0
$ws.Range("M49").ClearContents()

$ws.Range("H101").Value2 = 21885.666
$ws.Range("J101").Value2 = 21885.666
$ws.Range("L101").Value2 = 21885.666
$ws.Range("N101").Value2 = -28375.666

$ws.Range("H113").Value2 = 892.375
$ws.Range("I113").Value2 = 892.375
$ws.Range("J113").Value2 = 0
$ws.Range("K113").Value2 = 892.375
$ws.Range("L113").Value2 = 0
$ws.Range("M113").Value2 = 1277.625
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value2 = 150
$ws.Range("I25").Value2 = 150
$ws.Range("J25").Value2 = 0
$ws.Range("K25").Value2 = 450
$ws.Range("L25").Value2 = 0
$ws.Range("M25").Value2 = -281
$ws.Range("N25").ClearContents()

$ws.Range("H30").Value2 = 150
$ws.Range("I30").Value2 = 150
$ws.Range("J30").Value2 = 0
$ws.Range("K30").Value2 = 450
$ws.Range("L30").Value2 = 0
$ws.Range("M30").Value2 = -348
$ws.Range("N30").ClearContents()

$ws.Range("H55").Value2 = 965.1667
$ws.Range("I55").Value2 = 577.3333
$ws.Range("J55").Value2 = 1094.4445
$ws.Range("K55").Value2 = 1731.9999
$ws.Range("L55").Value2 = 3283.3335
$ws.Range("M55").Value2 = -1554.9999
$ws.Range("N55").Value2 = -3637.3335

$ws.Range("H131").Value2 = 899.9231
$ws.Range("J131").Value2 = 900
$ws.Range("L131").Value2 = 2700
$ws.Range("N131").Value2 = -12780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 14950
$ws.Range("I70").Value2 = 14950
$ws.Range("K70").Value2 = 14950
$ws.Range("M70").Value2 = -14680

$ws.Range("H73").Value2 = 14950
$ws.Range("I73").Value2 = 14950
$ws.Range("K73").Value2 = 14950
$ws.Range("M73").Value2 = -14014

$ws.Range("H80").Value2 = 1619
$ws.Range("I80").Value2 = 1750
$ws.Range("J80").Value2 = 1488
$ws.Range("K80").Value2 = 1750
$ws.Range("L80").Value2 = 1488
$ws.Range("M80").Value2 = -752
$ws.Range("N80").Value2 = -3484

$ws.Range("H83").Value2 = 1619
$ws.Range("I83").Value2 = 1750
$ws.Range("J83").Value2 = 1488
$ws.Range("K83").Value2 = 8750
$ws.Range("L83").Value2 = 7440
$ws.Range("M83").Value2 = -3758
$ws.Range("N83").Value2 = -17424

$ws.Range("H126").Value2 = 7247.375
$ws.Range("I126").Value2 = 7595.8
$ws.Range("K126").Value2 = 22787.4
$ws.Range("M126").Value2 = -20317.4

$ws.Range("H132").Value2 = 4095.6667
$ws.Range("I132").Value2 = 3649
$ws.Range("K132").Value2 = 10947
$ws.Range("M132").Value2 = -8417

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 1904.0769
$ws.Range("I22").Value2 = 1209.75
$ws.Range("J22").Value2 = 3015
$ws.Range("K22").Value2 = 1209.75
$ws.Range("L22").Value2 = 3015
$ws.Range("M22").Value2 = -914.75
$ws.Range("N22").Value2 = -3605

$ws.Range("H27").Value2 = 1904.0769
$ws.Range("I27").Value2 = 1209.75
$ws.Range("J27").Value2 = 3015
$ws.Range("K27").Value2 = 1209.75
$ws.Range("L27").Value2 = 3015
$ws.Range("M27").Value2 = -1102.75
$ws.Range("N27").Value2 = -3229

$ws.Range("H46").Value2 = 3013.7727
$ws.Range("I46").Value2 = 2253.7693
$ws.Range("K46").Value2 = 2253.7693
$ws.Range("M46").Value2 = -2065.7693

$ws.Range("H104").Value2 = 20370
$ws.Range("J104").Value2 = 20370
$ws.Range("L104").Value2 = 20370
$ws.Range("N104").Value2 = -27358

$ws.Range("H132").Value2 = 15665.5
$ws.Range("I132").Value2 = 15798.6
$ws.Range("K132").Value2 = 47395.8
$ws.Range("M132").Value2 = -44865.8

$ws.Range("H136").Value2 = 4199.25
$ws.Range("I136").Value2 = 3718.8
$ws.Range("K136").Value2 = 11156.4
$ws.Range("M136").Value2 = -8606.400000000001
